# --- Update database structure and entries ---
# Appends 93 new Kassite personal-name entries (rows 192-287, column A only)
# to the "Kassitische Analyse" sheet, matching existing header-row styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Abi-Rataš",
    "Abi-Rutaš",
    "Agum",
    "Alba-da",
    "Albadi-Saḫ",
    "Ališbi-Saḫ",
    "Ailzibu",
    "Alzibu-naša",
    "Alsa-duri",
    "Ašar-Saḫ",
    "Ani-kilamdi",
    "Ašriqa",
    "Ašriqu",
    "Badu-Zana",
    "Bi-Bugašu",
    "Biri-šuriyaš",
    "Bula-ḫali",
    "bula-nikir",
    "Bunabu",
    "Bunnie",
    "Bunni-ḫarpa",
    "Bunni-tura",
    "Burame",
    "Buramizi",
    "Bura-saḫ",
    "Bura-Sana",
    "Buraša-Marduk",
    "Bureya",
    "Burna-Adad",
    "Burnami-saḫ",
    "burna-Harbašu",
    "Burnabiḫe",
    "burnabu",
    "Burna-zini",
    "burni",
    "burni-durum",
    "burni-dura",
    "burni-mašḫu",
    "burpa-suqšu",
    "Burra-alban",
    "burra-laguda",
    "burraši-galdu",
    "burra-šikme",
    "burra-šuḫur",
    "burra-šuqamuna",
    "burra-turra",
    "burra-akmul",
    "burrutu",
    "duni-mašḫu",
    "ebi-rattaš",
    "gab-buriyaš",
    "gab-ḫarbe",
    "gab-satran",
    "gab-šuqamuna",
    "gaddaš",
    "gandaš",
    "gandi",
    "gerza-bura",
    "kirza-bura",
    "gunini-bugaš",
    "gunizar-bugaš",
    "gurba-saḫ",
    "gurpazaduya",
    "guzalzal-bugaš",
    "guzar-ali",
    "guzaru",
    "ḫamaš-šarri",
    "ḫamaš-šuk",
    "ḫāmaš-šugab",
    "ḫārbal-saḫ",
    "ḫarpaniwe",
    "ḫarba-šiḫu",
    "ḫašardu",
    "ḫašmar",
    "ḫašmar-galšu",
    "ḫašuar",
    "ḫirzi",
    "ḫudi-yazi",
    "ḫumurbi",
    "ḫumurti",
    "ḫuri-saḫ",
    "ḫušši-ḫarbe",
    "ḫuštiya",
    "Ibša-ḫalu",
    "ibša",
    "ili-šarigaš",
    "inza-ḫudak",
    "inzatum",
    "inzi-mašḫu",
    "inzite",
    "inzi-tešup",
    "inzu",
    "ibzu-ḫusieš",
    "inzu-kutir",
    "inzu-mena",
    "inzu-murudaš"
)

$startRow = 192
$lastRow = $startRow + $names.Length - 1

for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $names[$i]
}

# New rows reuse the same "Name" style (Arial, s=1) used for entries in A1:A191.
$ws.Range("A191").Copy()
$ws.Range("A" + $startRow + ":A" + $lastRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Refresh the frozen-pane view onto the newly added rows and restore selection/zoom.
$win = $excel.ActiveWindow
$win.Zoom = 291
$ws.Range("C280").Select() | Out-Null
